$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.128342032432556
$ws.Range("B1").Value = 2.658652544021606
$ws.Range("C1").Value = 5.767167091369629
$ws.Range("D1").Value = 2.089398145675659
$ws.Range("E1").Value = 1.199519753456116
